$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in / clear individual imputed values (rows 3-23) ---

# Row 3 (RM 8): D3 was blank -> now -14.2
$ws.Range("D3").Value = -14.2

# Row 4 (RM 9): E4 was -6.4 -> now blank
$ws.Range("E4").Value = $null

# Row 5 (RM 14): D5 was -14.4 -> now blank
$ws.Range("D5").Value = $null

# Row 9 (RM 42): E9 was blank -> now -6.8
$ws.Range("E9").Value = -6.8

# Row 10 (RM 52 a): E10 was blank -> now -6.1
$ws.Range("E10").Value = -6.1

# Row 17 (RM 116): E17 was -7.3 -> now blank
$ws.Range("E17").Value = $null

# Row 18 (RM 120): E18 was -8.5 -> now blank
$ws.Range("E18").Value = $null

# Row 21 (RM 135): D21 was blank -> now -14.3
$ws.Range("D21").Value = -14.3

# Row 23 (RM 140): D23 was -13.9 -> now blank
$ws.Range("D23").Value = $null

# --- Remove two data rows (RM 232 and SC 92), shifting everything up ---

# Row 26 currently holds "RM 232" - delete entirely
$ws.Rows.Item(26).Delete()

# After the delete above, "SC 92" (was row 28) is now at row 27 - delete entirely
$ws.Rows.Item(27).Delete()

# --- Fill in imputed value for SC 193 (now at row 32 after the two deletions) ---
$ws.Range("D32").Value = -14.7
